$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New summary / analysis block appended below the existing data
# (rows 74-76, columns E-H): a small pivot-style table computing
# counts / averages for "transportation" and "charging" tasks.
# ---------------------------------------------------------------

# Header row (row 74)
$ws.Range("E74").Value2 = "Task Type"
$ws.Range("F74").Value2 = "Count"
$ws.Range("G74").Value2 = "Total Time"
$ws.Range("H74").Value2 = "Waiting Time"

# Transportation summary row (row 75)
$ws.Range("E75").Value2 = "Transport"
$ws.Range("F75").Formula = '=COUNTIF(A2:A63,"transportation")'
$ws.Range("G75").Formula = '=AVERAGEIF(A2:A63,"transportation",F2:F63)'
$ws.Range("H75").Formula = '=AVERAGEIF(A2:A63,"transportation",H2:H63)'

# Charging summary row (row 76)
$ws.Range("E76").Value2 = "Charging"
$ws.Range("F76").Formula = '=COUNTIF(A2:A63,"charging")'
$ws.Range("G76").Formula = '=AVERAGEIF(A2:A63,"charging",F2:F63)'
$ws.Range("H76").Formula = '=AVERAGEIF(A2:A63,"charging",H2:H63)'

# ---------------------------------------------------------------
# Formatting. Build each compound style once on an out-of-the-way
# scratch cell, then paste just the formatting onto the real
# target range. Doing it this way (instead of setting Font/
# Interior/Borders/Alignment one by one directly on the target)
# avoids baking a separate intermediate cell style for every
# single property that gets touched along the way.
# ---------------------------------------------------------------

# Header style: bold font, light themed blue fill, thin box border, centered.
$scratch1 = $ws.Range("Z1")
$scratch1.Font.Bold = $true
$scratch1.Interior.Color = 15652797
$scratch1.Borders.LineStyle = 1
$scratch1.HorizontalAlignment = -4108

$header = $ws.Range("E74:H74")
$scratch1.Copy()
$header.PasteSpecial(-4122)  # xlPasteFormats
$scratch1.Clear() | Out-Null

# Body style: thin box border, centered (no fill/bold).
$scratch2 = $ws.Range("Z2")
$scratch2.Borders.LineStyle = 1
$scratch2.HorizontalAlignment = -4108

$body = $ws.Range("E75:H76")
$scratch2.Copy()
$body.PasteSpecial(-4122)  # xlPasteFormats
$scratch2.Clear() | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Update the view: select the newly added block, matching the
# author's on-screen state when saving.
# ---------------------------------------------------------------

$ws.Activate()
$ws.Range("A73:XFD75").Select()
